$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying NATMI computation was re-run with updated TPM data. The
# "Target cluster" = MuSCs rows are no longer produced by the new run, so
# those rows (originally rows 4, 7 and 10 - i.e. D column = "MuSCs") are
# removed entirely; the remaining rows keep their Sending/Target cluster
# pairing but get refreshed statistic columns (G:T).
# Delete bottom-up so earlier row numbers stay valid while deleting.
$ws.Range("A10:T10").EntireRow.Delete()
$ws.Range("A7:T7").EntireRow.Delete()
$ws.Range("A4:T4").EntireRow.Delete()

# Refresh the statistic columns (G through T) for the six remaining data
# rows (now rows 2-7) with the newly computed TPM-based values.
$ws.Range("G2").Value = 0.264496
$ws.Range("H2").Value = 0.793488
$ws.Range("I2").Value = 0.001006353962629067
$ws.Range("J2").Value = 0.001006353962629067
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.001937666666666667
$ws.Range("N2").Value = 0.005813
$ws.Range("O2").Value = 0.0230007399171451
$ws.Range("P2").Value = 0.02300073991714511
$ws.Range("Q2").Value = 0.0005125050826666666
$ws.Range("R2").Value = 0.004612545744
$ws.Range("S2").Value = 0.00002314688575901954
$ws.Range("T2").Value = 0.00002314688575901954

$ws.Range("G3").Value = 0.264496
$ws.Range("H3").Value = 0.793488
$ws.Range("I3").Value = 0.001006353962629067
$ws.Range("J3").Value = 0.001006353962629067
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.082306
$ws.Range("N3").Value = 0.246918
$ws.Range("O3").Value = 0.9769992600828549
$ws.Range("P3").Value = 0.976999260082855
$ws.Range("Q3").Value = 0.021769607776
$ws.Range("R3").Value = 0.195926469984
$ws.Range("S3").Value = 0.0009832070768700475
$ws.Range("T3").Value = 0.0009832070768700475

$ws.Range("G4").Value = 245.845932
$ws.Range("H4").Value = 737.537796
$ws.Range("I4").Value = 0.9353942133886188
$ws.Range("J4").Value = 0.9353942133886189
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.001937666666666667
$ws.Range("N4").Value = 0.005813
$ws.Range("O4").Value = 0.0230007399171451
$ws.Range("P4").Value = 0.02300073991714511
$ws.Range("Q4").Value = 0.4763674675719999
$ws.Range("R4").Value = 4.287307208147999
$ws.Range("S4").Value = 0.02151475902215415
$ws.Range("T4").Value = 0.02151475902215415

$ws.Range("G5").Value = 245.845932
$ws.Range("H5").Value = 737.537796
$ws.Range("I5").Value = 0.9353942133886188
$ws.Range("J5").Value = 0.9353942133886189
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.082306
$ws.Range("N5").Value = 0.246918
$ws.Range("O5").Value = 0.9769992600828549
$ws.Range("P5").Value = 0.976999260082855
$ws.Range("Q5").Value = 20.234595279192
$ws.Range("R5").Value = 182.111357512728
$ws.Range("S5").Value = 0.9138794543664647
$ws.Range("T5").Value = 0.9138794543664649

$ws.Range("G6").Value = 16.71558533333333
$ws.Range("H6").Value = 50.146756
$ws.Range("I6").Value = 0.06359943264875202
$ws.Range("J6").Value = 0.06359943264875202
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.001937666666666667
$ws.Range("N6").Value = 0.005813
$ws.Range("O6").Value = 0.0230007399171451
$ws.Range("P6").Value = 0.02300073991714511
$ws.Range("Q6").Value = 0.03238923251422222
$ws.Range("R6").Value = 0.2915030926279999
$ws.Range("S6").Value = 0.001462834009231932
$ws.Range("T6").Value = 0.001462834009231932

$ws.Range("G7").Value = 16.71558533333333
$ws.Range("H7").Value = 50.146756
$ws.Range("I7").Value = 0.06359943264875202
$ws.Range("J7").Value = 0.06359943264875202
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.082306
$ws.Range("N7").Value = 0.246918
$ws.Range("O7").Value = 0.9769992600828549
$ws.Range("P7").Value = 0.976999260082855
$ws.Range("Q7").Value = 1.375792966445333
$ws.Range("R7").Value = 12.382136698008
$ws.Range("S7").Value = 0.06213659863952009
$ws.Range("T7").Value = 0.0621365986395201
